$d = $word.ActiveDocument

# The paragraph currently reads (single run):
#   "...现在的解决方案是使用WSL2来启动一个Linux发行版"
# It needs to become (three runs, middle one new):
#   "...现在的解决方案是使用WSL2" + " + Docker" + "来启动一个Linux发行版"
#
# Insert " + Docker" right after "WSL2". We then nudge Font.Bold on/off for
# the freshly-inserted text (and again on the following "来启动一个Linux
# 发行版" text) so the run-coalescing pass on save keeps each piece in its
# own <w:r> instead of silently re-merging everything back into a single run.

$r = $d.Content
$found = $r.Find.Execute("WSL2", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'WSL2' in the document"
}
$r.Collapse(0)
$r.InsertAfter(" + Docker")
$r.Font.Bold = $true
$r.Font.Bold = $false

$r2 = $d.Content
$found2 = $r2.Find.Execute("来启动一个Linux发行版", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find '来启动一个Linux发行版' in the document"
}
$r2.Font.Bold = $true
$r2.Font.Bold = $false
